# Generate Report for Handoff
# Adds a new file (a8da677e-cb8e-42c1-ab73-2c51f3462c5b) handoff row to each
# sheet of the localization-status report, mirroring the existing
# 6ca3f72c-ceec-42b9-9f06-60dbef0651b1 row.

$wb = $excel.ActiveWorkbook

$newGuid = "a8da677e-cb8e-42c1-ab73-2c51f3462c5b"
$newHash = "b76129c4518c6779fb12fae9ab4249e431f49459"
$dateFmt = "yyyy-mm-dd HH:mm:ss"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/22ced7e827dc25355322fba6d0e1900e36b1f398/e2e/$newGuid.md"
$zhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bd95003d68d79094bd7d71000840622b2a6372b9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newGuid.$newHash.zh-cn.xlf"
$deUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a3694f19942011cf906b741b34b61934aec41fb1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newGuid.$newHash.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview": new row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "$newGuid.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-38-13 12:38:34"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdUrl, "", "", "$newGuid.md") | Out-Null
$wsOverview.Range("A3").Font.Underline = $true
$wsOverview.Range("A3").Font.Color = 15570276

# ---------------------------------------------------------------------------
# Sheet "zh-cn": new row 3
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = "$newGuid.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-13 12:38:31"
$wsZh.Range("E3").NumberFormat = $dateFmt
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl, "", "", "$newGuid.md") | Out-Null
$wsZh.Range("A3").Font.Underline = $true
$wsZh.Range("A3").Font.Color = 15570276

$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $mdUrl, "", "", ".md") | Out-Null
$wsZh.Range("B3").Font.Underline = $true
$wsZh.Range("B3").Font.Color = 15570276

$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhUrl, "", "", "$newGuid.$newHash.zh-cn.xlf") | Out-Null
$wsZh.Range("D3").Font.Underline = $true
$wsZh.Range("D3").Font.Color = 15570276

# ---------------------------------------------------------------------------
# Sheet "de-de": new row 3
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = "$newGuid.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-13 12:38:34"
$wsDe.Range("E3").NumberFormat = $dateFmt
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl, "", "", "$newGuid.md") | Out-Null
$wsDe.Range("A3").Font.Underline = $true
$wsDe.Range("A3").Font.Color = 15570276

$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $mdUrl, "", "", ".md") | Out-Null
$wsDe.Range("B3").Font.Underline = $true
$wsDe.Range("B3").Font.Color = 15570276

$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deUrl, "", "", "$newGuid.$newHash.de-de.xlf") | Out-Null
$wsDe.Range("D3").Font.Underline = $true
$wsDe.Range("D3").Font.Color = 15570276

Write-Host "Handoff report rows added."
